$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B34").Value = "62d0f019011e1e35afb4da08a53861dd"
$ws.Range("B89").Value = "3a425473b901d99eeb2f8f05d1a7a9da"
$ws.Range("B99").Value = "7332e19db9d80de1248db805e60f9312"
$ws.Range("B110").Value = "a0cab0e46f110ea81f706b2fc5953f20"
$ws.Range("B154").Value = "dadb7be999dbd63f806299bfafbc6261"
$ws.Range("B160").Value = "25264021f32130c246ff1dcdeec483d0"
$ws.Range("B162").Value = "d9cbdf45e33118bc240620a3976be092"
$ws.Range("B180").Value = "ae42a0af0e2092a422639ad4d71db265"
$ws.Range("B213").Value = "618db607106c4c865cbafcf8156b579a"
$ws.Range("B281").Value = "181895aa68478a8ce5e37e3a6123fdf6"
$ws.Range("B338").Value = "c16252edd9bbad81bece7e1e437aeca5"
$ws.Range("B468").Value = "e1e4b714dddf2e3deb6075c4d94ffcf9"
$ws.Range("B511").Value = "b3c0471f6ab03fe79ed3515cd46b22cc"
$ws.Range("B516").Value = "0f2b68cdf56bae47118f70f03e78d2f5"
$ws.Range("B524").Value = "3962d32114f3fb69ae6f12f86a119019"
$ws.Range("B535").Value = "320c9d5b1e38d46bf285d4beb72f820c"
$ws.Range("B545").Value = "6872b106d46507f66af37d33523f76f9"
$ws.Range("B559").Value = "a43aad2a42277be6fc85233bafe81f21"
$ws.Range("B565").Value = "2ba2af195a7150411e9edbf214040e44"
$ws.Range("B596").Value = "db79560a07b943a028661bf9ac58f8cf"
$ws.Range("B677").Value = "16b63d480f3d50d78a869c19ab998727"
$ws.Range("B678").Value = "654c1ba0472b17af82efd250300ae113"
$ws.Range("B741").Value = "1f9b18a75e7137204200fd2e581624f2"
$ws.Range("B780").Value = "7b32c2e2138ad20d6de90800ca768f42"
$ws.Range("B823").Value = "1240d1925d5bb6781d888325f1408e49"
$ws.Range("B827").Value = "18959c8132fbe58132b63e2ed262ede7"
$ws.Range("B828").Value = "683ad9d5a62eedccab952d06bed5a4f7"
$ws.Range("B837").Value = "c23d1d2e9e89bd032e026d27dfcc8827"
$ws.Range("B839").Value = "97010d418992034607b9ffb8ac4a8020"
